$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug in window type lengths test data: C4 should be 2, not 4
$ws.Range("C4").Value = 2

# Update the active cell selection
$ws.Range("O9").Select()
